# Minor updates to presentation
#
# 1. Slide 13 ("And yet more jargon!"): reword the WireMock bullet text.
# 2. Slide 19: nudge the "Group 9" callout group to a new position.
# 3. Slide 8 ("How can contract testing help?"): split a sentence into two
#    runs so that "before deployment" is bold.

$p = $ppt.ActivePresentation

# Helper: replace the first occurrence of $oldStr with $newStr inside a
# shape's TextFrame, using TextRange.Characters(start, length) (absolute,
# 1-based positions across the whole shape text -- this is reliable, unlike
# Paragraphs().Runs() whose Start/Length bookkeeping is unreliable once a
# paragraph has more than two runs).
function Replace-ShapeText($shapeTextRange, [string]$oldStr, [string]$newStr) {
    $full = $shapeTextRange.Text
    $idx = $full.IndexOf($oldStr)
    if ($idx -lt 0) {
        throw "Replace-ShapeText: substring not found: $oldStr"
    }
    $start = $idx + 1
    $len = $oldStr.Length
    $sub = $shapeTextRange.Characters($start, $len)
    $sub.Text = $newStr
}

# --- 1. Slide 13: WireMock wording -----------------------------------
$s13 = $p.Slides.Item(13)
$tf13 = $s13.Shapes.Item(2).TextFrame.TextRange

Replace-ShapeText $tf13 `
    "WireMock is a library for stubbing and mocking web services" `
    "WireMock is a server for stubbing and mocking web services.  It runs locally when running tests"

Replace-ShapeText $tf13 `
    "A micro-service sends requests to it, and WireMock sends replies" `
    "A unit test sends requests to it, and the WireMock server sends replies"

Replace-ShapeText $tf13 `
    "It is used by WireMock to work find a response which matches the request message" `
    "It is used by the WireMock server to find a response which matches the request message"

# --- 2. Slide 19: reposition the "Group 9" shape ----------------------
# The underlying engine stores Shape.Left/Top as single-precision floats,
# so a direct points<->EMU conversion can be off by 1 EMU. Bisect on the
# actual read-back value so the saved XML lands exactly on the target EMU.
function Set-ExactLeft($shape, [double]$targetEmu) {
    $lo = ($targetEmu - 2.0) / 914400.0 * 72.0
    $hi = ($targetEmu + 2.0) / 914400.0 * 72.0
    for ($iter = 0; $iter -lt 60; $iter++) {
        $mid = ($lo + $hi) / 2.0
        $shape.Left = $mid
        $readEmu = [math]::Round([double]$shape.Left / 72.0 * 914400.0)
        if ($readEmu -lt $targetEmu) { $lo = $mid } else { $hi = $mid }
    }
    $shape.Left = $hi
}
function Set-ExactTop($shape, [double]$targetEmu) {
    $lo = ($targetEmu - 2.0) / 914400.0 * 72.0
    $hi = ($targetEmu + 2.0) / 914400.0 * 72.0
    for ($iter = 0; $iter -lt 60; $iter++) {
        $mid = ($lo + $hi) / 2.0
        $shape.Top = $mid
        $readEmu = [math]::Round([double]$shape.Top / 72.0 * 914400.0)
        if ($readEmu -lt $targetEmu) { $lo = $mid } else { $hi = $mid }
    }
    $shape.Top = $hi
}

$s19 = $p.Slides.Item(19)
$grp = $s19.Shapes.Item(4)
Set-ExactLeft $grp 1267596
Set-ExactTop  $grp 4639580

# --- 3. Slide 8: bold "before deployment" ------------------------------
$s8 = $p.Slides.Item(8)
$tf8 = $s8.Shapes.Item(2).TextFrame.TextRange

$oldStr8 = "Contract testing should detect the problem from the previous slide before deployment"
$newFirstPart = "Contract testing should detect the problem shown in the previous slide "
Replace-ShapeText $tf8 $oldStr8 $newFirstPart

# Insert "before deployment" as its own run right after the first part.
$afterEdit = $tf8.Text
$idx = $afterEdit.IndexOf($newFirstPart)
$lastCharPos = $idx + 1 + $newFirstPart.Length - 1
$lastChar = $tf8.Characters($lastCharPos, 1)
$newRun = $lastChar.InsertAfter("before deployment")

# Bold only the newly-inserted text.
$afterInsert = $tf8.Text
$bdIdx = $afterInsert.IndexOf("before deployment")
$bdStart = $bdIdx + 1
$bdLen = "before deployment".Length
$bdRange = $tf8.Characters($bdStart, $bdLen)
$bdRange.Font.Bold = $true

# InsertAfter pulled the preceding space into its own (non-bold) run;
# fold it back into the first run so we end up with exactly two runs:
# "...previous slide " (regular) + "before deployment" (bold).
$spaceRange = $tf8.Characters($bdStart - 1, 1)
$spaceRange.Text = ""

$full4 = $tf8.Text
$r1Idx = $full4.IndexOf("Contract testing should detect the problem shown in the previous slide")
$r1Start = $r1Idx + 1
$r1Len = "Contract testing should detect the problem shown in the previous slide".Length
$r1Range = $tf8.Characters($r1Start, $r1Len)
$r1Range.Text = $newFirstPart
